# Apply the LOM3260 syllabus content update:
#  - Fix the (Portuguese) "Objetivos:" row content.
#  - Split the single "docentes responsaveis" row into two rows (two lecturers).
#  - Add a Portuguese "short syllabus" row and a Portuguese full "Programa" row.
#  - Correct the Method / Criterio / Norma de recuperacao / Bibliografia rows,
#    which had been shifted by one in the source file, and add the
#    "Bibliografia" text itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the Objetivos (Portuguese) text in row 10 -----------------------
$objetivosPt = "Fornecer ao aluno uma introdução à computação científica moderna, usando a linguagem Python e suas bibliotecas numéricas e gráficas mais populares: numpy, scipy, matplotlib e pandas. Ao final do curso, o aluno estará capacitado a desenvolver programas complexos, de pequeno e médio porte para solucionar problemas de engenharia que envolvam processamento numérico de grandes conjuntos de dados e correlacionar variáveis usando métodos numéricos."
$ws.Range("B10").Value = $objetivosPt
$ws.Range("C10").Value = $objetivosPt

# --- 2. Insert two new rows after row 12 (Docentes responsáveis:) ----------
# These become the rows for the two listed instructors.
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

# Row 13: first instructor
$emerson = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("B13").Value = $emerson
$ws.Range("C13").Value = $emerson

# Row 14: second instructor
$luiz = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("B14").Value = $luiz
$ws.Range("C14").Value = $luiz

# --- 3. Row 15 "Programa resumido:" now needs the Portuguese short syllabus -
$resumidoPt = "Introdução à programação em Python; palavras-chave em Python; rotinas e funções; classes; numpy e o conceito de slicing e indexing de arrays; revisão de métodos numéricos usando scipy; geração de gráficos e animações com a biblioteca matplotlib; criação de interfaces gráficas com o usuário usando matplotlib.widgets"
$ws.Range("B15").Value = $resumidoPt
$ws.Range("C15").Value = $resumidoPt

# --- 4. Row 17 "Programa:" needs the Portuguese full syllabus --------------
$programaPt = [char]0x2022 + " Introdução à programação em Python " + [char]0x2022 + " Instalação de uma distribuição Python em Windows e Linux " + [char]0x2022 + " Formatação de arquivos em Python " + [char]0x2022 + " Estruturas condicionais " + [char]0x2022 + " Laços de repetição de comandos " + [char]0x2022 + " Outras palavras-chaves e métodos " + [char]0x2022 + " Rotinas e funções " + [char]0x2022 + " Códigos multifonte e bibliotecas pessoais " + [char]0x2022 + " Bibliotecas numéricas e gráficas: numpy, scipy e matplotlib " + [char]0x2022 + " Programação orientada a objeto: classes " + [char]0x2022 + " Conceito de objetos e instâncias " + [char]0x2022 + " Classes e subclasses" + [char]0x2022 + " " + [char]0x201C + "Arrays" + [char]0x201D + " em numpy " + [char]0x2022 + " O conceito de array em numpy " + [char]0x2022 + " " + [char]0x201C + "Slicing" + [char]0x201D + " e indexação " + [char]0x2022 + " Trabalhando com arquivos (entrada e saída) " + [char]0x2022 + " Gráficos em matplotlib " + [char]0x2022 + " A biblioteca matplotlib.pyplot e gráficos em 2D e 3D " + [char]0x2022 + " A biblioteca matplotlib.animation para criar gráficos animados. " + [char]0x2022 + " Interfaces gráficas com o usuário (Graphical User Interface, GUI) " + [char]0x2022 + " Interfaces simples com a biblioteca matplotlib.widgets."
$ws.Range("B17").Value = $programaPt
$ws.Range("C17").Value = $programaPt

# --- 5. Row 20 "Método:" text -----------------------------------------------
$metodo = "Aulas expositivas e em laboratório computacional, trabalhos e exercícios comentados."
$ws.Range("B20").Value = $metodo
$ws.Range("C20").Value = $metodo

# --- 6. Row 21 "Critério:" text ---------------------------------------------
$criterio = "Média aritmética de exercícios e trabalhos propostos ao longo do curso e uma apresentação final de projeto."
$ws.Range("B21").Value = $criterio
$ws.Range("C21").Value = $criterio

# --- 7. Row 22 "Norma de recuperação:" text --------------------------------
$norma = "Não haverá exame de recuperação."
$ws.Range("B22").Value = $norma
$ws.Range("C22").Value = $norma

# --- 8. Row 23 "Bibliografia:" text -----------------------------------------
$bib = "Lambert, K. A. Fundamentos de Python: estruturas de dados. Cengage, 2ed, 2022.Nilo Ney Coutinho Menezes. Introdução à Programação com Python: Algoritmos e Lógica de Programação Para Iniciantes, 3a ed, 2019.Ramalho, L. Python Fluente. O" + [char]0x2019 + "Reilly-Novatec, 2015Downey, A. B. Pense em Python. O" + [char]0x2019 + "Reilly-Novatec, 2016.STEWART, J. M. Python for scientists. Cambridge University Press, 2014.TELLES, M. Python Power, Boston: Thomson Course Technology PTR, 2008.LUTZ, Mark. Programming Python, 3a ed, Sebastopol, CA: O" + [char]0x2019 + "Reilly Media, 2006.MCGREGGOR, D. M. Mastering matplotlib. Birmingham, UK: Packt Publishing, 2015."
$ws.Range("B23").Value = $bib
$ws.Range("C23").Value = $bib

# --- 9. Row height bookkeeping ----------------------------------------------
# Rows 13 and 14 (new instructor rows) keep the default row height (no custom
# height), matching the other simple "label: value" rows such as row 12 -- so
# nothing to set for them; newly inserted rows start out at default height.

# Rows 15/16 (short syllabus, PT/EN) -> 60pt tall.
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 60

# Rows 17/18 (full syllabus, PT/EN) -> 120pt tall.
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 120

# Rows 20/21/22 (Método/Critério/Norma) -> 60pt tall.
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 60

# Row 23 (Bibliografia) -> 120pt tall.
$ws.Rows.Item(23).RowHeight = 120

